$d = $word.ActiveDocument

# The three <wp:drawing> logos (2x PearsonLogo footers + 1x BTec_Logo header)
# were re-exported with swapped image part "display names" (the `name`
# attribute on both the drawing's <wp:docPr> and the picture's <pic:cNvPr>).
# Neither element is reachable as a distinct, settable COM property pair on
# InlineShape (InlineShape.Name only ever touches <wp:docPr name>), so patch
# the underlying package XML directly via Document.WordOpenXML, which
# round-trips losslessly through this document's parts.

$xml = $d.WordOpenXML

$xml = $xml.Replace(
    '<wp:docPr descr="BTec_Logo-Orange" id="1" name="image1.jpg"/>',
    '<wp:docPr descr="BTec_Logo-Orange" id="1" name="image2.jpg"/>')
$xml = $xml.Replace(
    '<pic:cNvPr descr="BTec_Logo-Orange" id="0" name="image1.jpg"/>',
    '<pic:cNvPr descr="BTec_Logo-Orange" id="0" name="image2.jpg"/>')

$pearsonDesc = 'Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png'

$xml = $xml.Replace(
    '<wp:docPr descr="' + $pearsonDesc + '" id="2" name="image2.png"/>',
    '<wp:docPr descr="' + $pearsonDesc + '" id="2" name="image1.png"/>')
$xml = $xml.Replace(
    '<wp:docPr descr="' + $pearsonDesc + '" id="3" name="image2.png"/>',
    '<wp:docPr descr="' + $pearsonDesc + '" id="3" name="image1.png"/>')
$xml = $xml.Replace(
    '<pic:cNvPr descr="' + $pearsonDesc + '" id="0" name="image2.png"/>',
    '<pic:cNvPr descr="' + $pearsonDesc + '" id="0" name="image1.png"/>')

$d.WordOpenXML = $xml

Write-Output "done"
